# Apply the "Added Commodore 64 section" update.
#
# The workbook recounts lines/words for several source-code trees. The
# underlying per-file counts live as literal (non-formula) numbers on the
# "linecount.txt" and "wordcount.txt" sheets; the "Totals" sheet only
# contains SUMIFS/SUM formulas over that data, so it recalculates itself
# once the raw numbers below are updated.

$wb = $excel.ActiveWorkbook

# ---- linecount.txt sheet -------------------------------------------------
$lc = $wb.Worksheets.Item("linecount.txt")

$lc.Range("C38").Value = 2002363

$lc.Range("C50").Value = 10107

$lc.Range("A53").Value = 50776
$lc.Range("B53").Value = 283174
$lc.Range("C53").Value = 1982569

$lc.Range("B57").Value = 320
$lc.Range("C57").Value = 2791

$lc.Range("B58").Value = 665
$lc.Range("C58").Value = 5092

$lc.Range("A59").Value = 1536
$lc.Range("B59").Value = 4655
$lc.Range("C59").Value = 33475

$lc.Range("A60").Value = 1647
$lc.Range("B60").Value = 9421
$lc.Range("C60").Value = 65662

$lc.Range("A65").Value = 42473
$lc.Range("B65").Value = 216139
$lc.Range("C65").Value = 1511272

$lc.Range("A69").Value = 9688
$lc.Range("B69").Value = 46430
$lc.Range("C69").Value = 327338

$lc.Range("A72").Value = 48441
$lc.Range("B72").Value = 271141
$lc.Range("C72").Value = 1875698

$lc.Range("C112").Value = 140274

$lc.Range("B118").Value = 27663
$lc.Range("C118").Value = 214850

# ---- wordcount.txt sheet --------------------------------------------------
$wc = $wb.Worksheets.Item("wordcount.txt")

$wc.Range("A53").Value = 215405
$wc.Range("A57").Value = 195
$wc.Range("A58").Value = 237
$wc.Range("A59").Value = 2594
$wc.Range("A60").Value = 5984
$wc.Range("A65").Value = 160521
$wc.Range("A72").Value = 205992
$wc.Range("A118").Value = 19922

# Make sure every dependent formula (Totals sheet) is refreshed before save.
$excel.CalculateFullRebuild()
$excel.Calculate()
